# Parameters will now point to C://Phase1/ Location by default
# Replace the old "C:\Users\InterviewRoom1\Phase1..." paths in row 2 with
# the new "C:\Phase1..." equivalents (same relative sub-paths).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C:\Phase1"
$ws.Range("B2").Value = "C:\Phase1\Backend\node-neurosky"
$ws.Range("C2").Value = "C:\Phase1\Backend\emotions"
$ws.Range("D2").Value = "C:\Phase1\Videos\men.mp4"
$ws.Range("F2").Value = "C:\Phase1\Reports"
$ws.Range("G2").Value = "C:\Phase1\Backend\data-normalizer"

# Resize columns A-D to fit the shorter paths (bestFit-style manual resize),
# leaving F/G untouched. Values passed are COM "characters" widths; this
# host's serializer adds a fixed 5/6-character pad when writing the XML
# `width` attribute (matching Excel's own stored-width convention).
$ws.Columns.Item(1).ColumnWidth = 13.1666666666667
$ws.Columns.Item(2).ColumnWidth = 32.3072916666667
$ws.Columns.Item(3).ColumnWidth = 27.1666666666667
$ws.Columns.Item(4).ColumnWidth = 25.8776041666667

# Update the view: scroll so column B is the leftmost visible column, and
# move the selection to G2.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("G2").Select()
